$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / styles) of the last existing data
# row (row 5) down into the two new rows before filling in their values, so
# the new rows pick up the same cell styles (date style, accounting styles)
# as the rest of the table instead of the workbook's default style.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A5:F5").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

# Row 6: 四方坪站 station, 2025-12-03 (serial date 45994)
$ws.Cells.Item(6, 1).Value = 45994
$ws.Cells.Item(6, 2).Value = "四方坪站"
$ws.Cells.Item(6, 3).Value = 8653.2099999999991
$ws.Cells.Item(6, 4).Value = 7601.29
$ws.Cells.Item(6, 5).Value = 2904.55
$ws.Cells.Item(6, 6).Value = 373

# Row 7: 高岭站 station, 2025-12-03 (serial date 45994)
$ws.Cells.Item(7, 1).Value = 45994
$ws.Cells.Item(7, 2).Value = "高岭站"
$ws.Cells.Item(7, 3).Value = 4352.49
$ws.Cells.Item(7, 4).Value = 3655.02
$ws.Cells.Item(7, 5).Value = 1190.8399999999999
$ws.Cells.Item(7, 6).Value = 165

# Match the author's final selection state in the saved workbook.
$ws.Range("G11").Select()
